$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "defaultModel" header in B1 to "model5" (unique model name)
$ws.Range("B1").Value = "model5"

# Add new model columns E, F, G with their header names
$ws.Range("E1").Value = "model7"
$ws.Range("F1").Value = "model9"
$ws.Range("G1").Value = "model100"

# Fill in the feature-support matrix (1 = supported) for the new columns
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("E6").Value = 1

$ws.Range("E7").Value = 1

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1

$ws.Range("G10").Value = 1

$ws.Range("G11").Value = 1
